$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated fight-card data (rows 2..27): First, Last, Prediction
$data = @(
    @("Brendan", "Allen", 0),
    @("Paul", "Craig", 0),
    @("Michael", "Morales", 1),
    @("Jake", "Matthews", 0),
    @("Chase", "Hooper", 1),
    @("Jordan", "Leavitt", 1),
    @("Payton", "Talbott", 1),
    @("Nick", "Aguirre", 1),
    @("Amanda", "Ribas", 0),
    @("Luana", "Pinheiro", 0),
    @("Myktybek", "Orolbai", 1),
    @("Uros", "Medic", 1),
    @("Joanderson", "Brito", 0),
    @("Jonathan", "Pearce", 0),
    @("Jose", "Johnson", 0),
    @("Chad", "Anheliger", 0),
    @("Christian", "Leroy Duncan", 1),
    @("Denis", "Tiuliulin", 1),
    @("Mick", "Parkin", 0),
    @("Caio", "Machado", 0),
    @("Jeka", "Saragih", 1),
    @("Lucas", "Alexander", 0),
    @("Ailin", "Perez", 0),
    @("Lucie", "Pudilova", 0),
    @("Trey", "Ogden", 1),
    @("Nikolas", "Motta", 1)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $row++
}

# Fill in the index column (A) for the two newly added rows
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(27, 1).Value = 25

# Copy A25's formatting (bold, centered, bordered) down onto the new rows
$ws.Range("A25").Copy()
$ws.Range("A26:A27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
